$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: new "person-1.spec" test header row (values entered first so the
# shared-string table gets populated in source order: D9, E9, then E11, E12, E10)
$ws.Range("D9").Value = "person-1.spec"
$ws.Range("E9").Value = "persons"
$ws.Range("F9").Value = 23
$ws.Range("G9").Value = 5
$ws.Range("H9").Value = 33
$ws.Range("I9").Value = 4
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = 11
$ws.Range("L9").Value = 15

# Row 11 ("overview (person)")
$ws.Range("E11").Value = "overview (person)"
$ws.Range("F11").Value = 38
$ws.Range("G11").Value = 12
$ws.Range("H11").Value = 49

# Row 12 ("base (person)")
$ws.Range("E12").Value = "base (person)"
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0

# Row 10 ("person")
$ws.Range("E10").Value = "person"
$ws.Range("F10").Value = 100
$ws.Range("G10").Value = 100
$ws.Range("H10").Value = 100

# Border/formatting tweaks on rows 10-11 matching the underlying edit
# (bottom border cleared on D10/E10/F10/H10; all borders cleared on
# G10/F11/G11 and the I:M detail columns for rows 10-11)
$ws.Range("D10").Borders(9).LineStyle = -4142
$ws.Range("E10").Borders(9).LineStyle = -4142
$ws.Range("F10").Borders(9).LineStyle = -4142
$ws.Range("H10").Borders(9).LineStyle = -4142
$ws.Range("G10").Borders.LineStyle = -4142
$ws.Range("I10:M10").Borders.LineStyle = -4142
$ws.Range("F11").Borders.LineStyle = -4142
$ws.Range("G11").Borders.LineStyle = -4142
$ws.Range("I11:M11").Borders.LineStyle = -4142

# Sheet view / selection changes
$ws.Application.Goto($ws.Range("A4"), $false)
$ws.Range("M16").Select()
